# Auto-generated from the cryptos.xlsx OOXML diff.
# Column D (Price) cells are plain numeric-looking text in the source sheet
# (e.g. "0.536", "51.185.61"); a leading apostrophe forces Excel to keep them
# as text instead of auto-converting to a Number, matching the original t="inlineStr" cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'51.185.61"
$ws.Range("E2").Value = "  +0.25%  "
$ws.Range("D3").Value = "'2.954.18"
$ws.Range("E3").Value = "  +0.27%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "'375.55"
$ws.Range("E5").Value = "  -0.71%  "
$ws.Range("D6").Value = "'102.66"
$ws.Range("E6").Value = "  -1.71%  "
$ws.Range("D7").Value = "'0.536"
$ws.Range("E7").Value = "  -1.03%  "
$ws.Range("E8").Value = "  +0.12%  "
$ws.Range("D9").Value = "'0.585"
$ws.Range("E9").Value = "  -1.15%  "
$ws.Range("D10").Value = "'36.51"
$ws.Range("E10").Value = "  -1.26%  "
$ws.Range("E11").Value = "  -0.24%  "
$ws.Range("D12").Value = "'0.0839"
$ws.Range("E12").Value = "  +0.07%  "
$ws.Range("D13").Value = "'3.409.30"
$ws.Range("E13").Value = "  +0.01%  "
$ws.Range("D14").Value = "'17.96"
$ws.Range("E14").Value = "  -2.23%  "
$ws.Range("D15").Value = "'7.39"
$ws.Range("E15").Value = "  -1.13%  "
$ws.Range("D16").Value = "'2.958.90"
$ws.Range("E16").Value = "  +0.88%  "
$ws.Range("D17").Value = "'0.986"
$ws.Range("E17").Value = "  +2.53%  "
$ws.Range("D18").Value = "'51.124.34"
$ws.Range("E18").Value = "  +0.17%  "
$ws.Range("B19").Value = "Uniswap"
$ws.Range("C19").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D19").Value = "'8.03"
$ws.Range("E19").Value = "  +9.31%  "
$ws.Range("B20").Value = "ImmutableX"
$ws.Range("C20").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D20").Value = "'3.16"
$ws.Range("E20").Value = "  -4.79%  "
$ws.Range("D21").Value = "'12.69"
$ws.Range("E21").Value = "  -1.24%  "
$ws.Range("D22").Value = "'0.0₃0958"
$ws.Range("E22").Value = "  -0.03%  "
$ws.Range("D23").Value = "'264.79"
$ws.Range("E23").Value = "  +1.52%  "
$ws.Range("D24").Value = "'68.42"
$ws.Range("E24").Value = "  -1.32%  "
$ws.Range("D25").Value = "'2.90"
$ws.Range("E25").Value = "  +3.25%  "
$ws.Range("D26").Value = "'8.52"
$ws.Range("E26").Value = "  +12.22%  "
$ws.Range("D27").Value = "'8.05"
$ws.Range("E27").Value = "  +10.56%  "
$ws.Range("B28").Value = "Kaspa"
$ws.Range("C28").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D28").Value = "'0.170"
$ws.Range("E28").Value = "  +0.46%  "
$ws.Range("B29").Value = "Hedera"
$ws.Range("C29").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D29").Value = "'0.114"
$ws.Range("E29").Value = "  +1.20%  "
$ws.Range("E30").Value = "  +0.03%  "
$ws.Range("D31").Value = "'25.70"
$ws.Range("E31").Value = "  -0.36%  "
$ws.Range("D32").Value = "'9.90"
$ws.Range("E32").Value = "  +0.79%  "
$ws.Range("D33").Value = "'50.90"
$ws.Range("E33").Value = "  +0.09%  "
$ws.Range("D34").Value = "'33.87"
$ws.Range("E34").Value = "  -1.97%  "
$ws.Range("D35").Value = "'0.0449"
$ws.Range("E35").Value = "  +0.54%  "
$ws.Range("D36").Value = "'2.03"
$ws.Range("E36").Value = "  -2.73%  "
$ws.Range("E37").Value = "  -0.30%  "
$ws.Range("D38").Value = "'3.00"
$ws.Range("E38").Value = "  -1.68%  "
$ws.Range("D39").Value = "'2.58"
$ws.Range("E39").Value = "  +0.39%  "
$ws.Range("E40").Value = "  +0.21%  "
$ws.Range("D41").Value = "'16.50"
$ws.Range("E41").Value = "  -4.03%  "
$ws.Range("D42").Value = "'1.80"
$ws.Range("E42").Value = "  -2.05%  "
$ws.Range("D43").Value = "'120.81"
$ws.Range("E43").Value = "  -1.17%  "
$ws.Range("D44").Value = "'0.289"
$ws.Range("E44").Value = "  +0.84%  "
$ws.Range("D45").Value = "'21.10"
$ws.Range("E45").Value = "  -3.86%  "
$ws.Range("D46").Value = "'2.03"
$ws.Range("E46").Value = "  -1.75%  "
$ws.Range("D47").Value = "'3.28"
$ws.Range("E47").Value = "  +2.39%  "
$ws.Range("E48").Value = "  -3.74%  "
$ws.Range("D49").Value = "'1.981.29"
$ws.Range("E49").Value = "  -2.43%  "
$ws.Range("D50").Value = "'0.0347"
$ws.Range("E50").Value = "  +1.19%  "
$ws.Range("D51").Value = "'5.07"
$ws.Range("E51").Value = "  -0.25%  "
